# Auto-generated edit script: applies numeric value updates to the
# Adamantoise_Profits profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each sheet holds static (non-formula) crafting-profit data in columns H:N;
# this script patches individual cell values, clears cells that should no
# longer be present, and adds cells that are newly populated.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 215.57143
$ws.Range("I8").Value = 215.57143
$ws.Range("K8").Value = 646.71429
$ws.Range("M8").Value = -507.71429
$ws.Range("H62").Value = 9681.5
$ws.Range("I62").Value = 7916.5
$ws.Range("K62").Value = 7916.5
$ws.Range("M62").Value = -7292.5
$ws.Range("H65").Value = 9681.5
$ws.Range("I65").Value = 7916.5
$ws.Range("K65").Value = 39582.5
$ws.Range("M65").Value = -36462.5
$ws.Range("H80").Value = 38466320
$ws.Range("I80").Value = 90909496
$ws.Range("J80").Value = 7991.2
$ws.Range("K80").Value = 272728488
$ws.Range("L80").Value = 23973.6
$ws.Range("M80").Value = -272727490
$ws.Range("N80").Value = -25969.6
$ws.Range("H83").Value = 38466320
$ws.Range("I83").Value = 90909496
$ws.Range("J83").Value = 7991.2
$ws.Range("K83").Value = 818185464
$ws.Range("L83").Value = 71920.8
$ws.Range("M83").Value = -818180472
$ws.Range("N83").Value = -81904.8
$ws.Range("H101").Value = 472.2143
$ws.Range("I101").Value = 478.75
$ws.Range("K101").Value = 1436.25
$ws.Range("M101").Value = 185.75
$ws.Range("H138").Value = 2415.4092
$ws.Range("I138").Value = 1846.3572
$ws.Range("J138").Value = 2834.7104
$ws.Range("K138").Value = 5539.071599999999
$ws.Range("L138").Value = 8504.1312
$ws.Range("M138").Value = -399.0715999999993
$ws.Range("N138").Value = -18784.1312

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2120.0435
$ws.Range("J2").Value = 3501.8572
$ws.Range("L2").Value = 3501.8572
$ws.Range("N2").Value = -3727.8572
$ws.Range("H4").Value = 202
$ws.Range("I4").Value = 202
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 202
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -86
$ws.Range("N4").ClearContents()
$ws.Range("H32").Value = 7442742
$ws.Range("I32").Value = 3789003
$ws.Range("J32").Value = 20839786
$ws.Range("K32").Value = 3789003
$ws.Range("L32").Value = 20839786
$ws.Range("M32").Value = -3788716
$ws.Range("N32").Value = -20840360
$ws.Range("H102").Value = 1185.2858
$ws.Range("I102").Value = 1088.5555
$ws.Range("J102").Value = 1359.4
$ws.Range("K102").Value = 1088.5555
$ws.Range("L102").Value = 1359.4
$ws.Range("M102").Value = 533.4445000000001
$ws.Range("N102").Value = -4603.4
$ws.Range("H116").Value = 2120.0435
$ws.Range("J116").Value = 3501.8572
$ws.Range("L116").Value = 3501.8572
$ws.Range("N116").Value = -8089.8572
$ws.Range("H122").Value = 4520.97
$ws.Range("I122").Value = 3588.2327
$ws.Range("K122").Value = 10764.6981
$ws.Range("M122").Value = -8314.6981
$ws.Range("H128").Value = 149980
$ws.Range("J128").Value = 149980
$ws.Range("L128").Value = 149980
$ws.Range("N128").Value = -159940
$ws.Range("H132").Value = 2441.1777
$ws.Range("I132").Value = 1804.9706
$ws.Range("K132").Value = 5414.9118
$ws.Range("M132").Value = -2884.9118

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2120.0435
$ws.Range("J3").Value = 3501.8572
$ws.Range("L3").Value = 3501.8572
$ws.Range("N3").Value = -3729.8572
$ws.Range("H94").Value = 1363.6316
$ws.Range("I94").Value = 896.73334
$ws.Range("J94").Value = 3114.5
$ws.Range("K94").Value = 896.73334
$ws.Range("L94").Value = 3114.5
$ws.Range("M94").Value = -445.73334
$ws.Range("N94").Value = -4016.5
$ws.Range("H132").Value = 114775
$ws.Range("J132").Value = 114775
$ws.Range("L132").Value = 114775
$ws.Range("N132").Value = -124895

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 280
$ws.Range("J22").Value = 350
$ws.Range("L22").Value = 350
$ws.Range("N22").Value = -1050
$ws.Range("H75").Value = 102325.664
$ws.Range("J75").Value = 102325.664
$ws.Range("L75").Value = 102325.664
$ws.Range("N75").Value = -104321.664
$ws.Range("H78").Value = 102325.664
$ws.Range("J78").Value = 102325.664
$ws.Range("L78").Value = 306976.992
$ws.Range("N78").Value = -316960.992
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H100").Value = 110890
$ws.Range("J100").Value = 110890
$ws.Range("L100").Value = 110890
$ws.Range("N100").Value = -113054
$ws.Range("H105").Value = 2538
$ws.Range("I105").Value = 2245.6
$ws.Range("K105").Value = 2245.6
$ws.Range("M105").Value = -498.5999999999999
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H141").Value = 1145816
$ws.Range("I141").Value = 59996
$ws.Range("K141").Value = 59996
$ws.Range("M141").Value = -54816

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 420
$ws.Range("I13").Value = 420
$ws.Range("K13").Value = 1260
$ws.Range("M13").Value = -1092
$ws.Range("H18").Value = 4850.8335
$ws.Range("I18").Value = 109
$ws.Range("K18").Value = 327
$ws.Range("M18").Value = -158
$ws.Range("H131").Value = 1757.2667
$ws.Range("J131").Value = 1858.1
$ws.Range("L131").Value = 5574.299999999999
$ws.Range("N131").Value = -15654.3
$ws.Range("H134").Value = 1866.6666
$ws.Range("I134").Value = 1866.6666
$ws.Range("K134").Value = 5599.9998
$ws.Range("M134").Value = -529.9997999999996
$ws.Range("H140").Value = 50002124
$ws.Range("I140").Value = 50002124
$ws.Range("K140").Value = 150006372
$ws.Range("M140").Value = -150001192

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 39999
$ws.Range("J39").Value = 39999
$ws.Range("L39").Value = 39999
$ws.Range("N39").Value = -41063
$ws.Range("H43").Value = 5198
$ws.Range("J43").Value = 7200
$ws.Range("L43").Value = 7200
$ws.Range("N43").Value = -7502
$ws.Range("H46").Value = 30500
$ws.Range("J46").Value = 39333.332
$ws.Range("L46").Value = 39333.332
$ws.Range("N46").Value = -39645.332
$ws.Range("H80").Value = 2047
$ws.Range("J80").Value = 2775.8
$ws.Range("L80").Value = 2775.8
$ws.Range("N80").Value = -4771.8
$ws.Range("H83").Value = 2047
$ws.Range("J83").Value = 2775.8
$ws.Range("L83").Value = 13879
$ws.Range("N83").Value = -23863
$ws.Range("H103").Value = 66729.664
$ws.Range("J103").Value = 66729.664
$ws.Range("L103").Value = 66729.664
$ws.Range("N103").Value = -69073.664
$ws.Range("H132").Value = 3869.5264
$ws.Range("I132").Value = 3010.3333
$ws.Range("J132").Value = 5342.4287
$ws.Range("K132").Value = 9030.999899999999
$ws.Range("L132").Value = 16027.2861
$ws.Range("M132").Value = -6500.999899999999
$ws.Range("N132").Value = -21087.2861
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2327.0625
$ws.Range("I22").Value = 2131.7273
$ws.Range("J22").Value = 2756.8
$ws.Range("K22").Value = 2131.7273
$ws.Range("L22").Value = 2756.8
$ws.Range("M22").Value = -1836.7273
$ws.Range("N22").Value = -3346.8
$ws.Range("H27").Value = 2327.0625
$ws.Range("I27").Value = 2131.7273
$ws.Range("J27").Value = 2756.8
$ws.Range("K27").Value = 2131.7273
$ws.Range("L27").Value = 2756.8
$ws.Range("M27").Value = -2024.7273
$ws.Range("N27").Value = -2970.8
$ws.Range("H134").Value = 43999
$ws.Range("J134").Value = 43999
$ws.Range("L134").Value = 43999
$ws.Range("N134").Value = -54139

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 556257.75
$ws.Range("I4").Value = 790
$ws.Range("J4").Value = 5000000
$ws.Range("K4").Value = 790
$ws.Range("L4").Value = 5000000
$ws.Range("M4").Value = -677
$ws.Range("N4").Value = -5000226
$ws.Range("H81").Value = 7048.1665
$ws.Range("I81").Value = 4096.3335
$ws.Range("K81").Value = 8192.666999999999
$ws.Range("M81").Value = -7131.666999999999
$ws.Range("H84").Value = 7048.1665
$ws.Range("I84").Value = 4096.3335
$ws.Range("K84").Value = 40963.335
$ws.Range("M84").Value = -35659.335
$ws.Range("H122").Value = 3714.6365
$ws.Range("I122").Value = 3179.1667
$ws.Range("J122").Value = 5142.5557
$ws.Range("K122").Value = 9537.500100000001
$ws.Range("L122").Value = 15427.6671
$ws.Range("M122").Value = -7087.500100000001
$ws.Range("N122").Value = -20327.6671
$ws.Range("H130").Value = 119994
$ws.Range("J130").Value = 119994
$ws.Range("L130").Value = 119994
$ws.Range("N130").Value = -130034
$ws.Range("H132").Value = 2698.9023
$ws.Range("I132").Value = 2190.743
$ws.Range("K132").Value = 6572.228999999999
$ws.Range("M132").Value = -4042.228999999999
$ws.Range("H136").Value = 2621.6316
$ws.Range("I136").Value = 2063.3928
$ws.Range("K136").Value = 6190.178400000001
$ws.Range("M136").Value = -3640.178400000001
